$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.211.65"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "3.049.12"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.74"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.34"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.27"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.371"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "3.576.17"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.28"
$ws.Range("E14").Value = "  -3.93%  "
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "57.325.77"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.12"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.041.26"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.03"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.22"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.497"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.76"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "0.0₃0894"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.28"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.72"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.28"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("B34").Value = "EnergySwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.84"
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "3.089.24"
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.90"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.650"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "2.249.91"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("E45").Value = "  +6.95%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.21"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.83"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "260.75"
$ws.Range("E50").Value = "  +13.46%  "
$ws.Range("E51").Value = "  +1.68%  "
